$wb = $excel.ActiveWorkbook

$wsCalorie  = $wb.Worksheets.Item("calorie_df")
$wsWeight   = $wb.Worksheets.Item("weight_df")
$wsExercise = $wb.Worksheets.Item("exercise_df")

# ---------------------------------------------------------------------------
# 1. calorie_df: add "Weight" and "Steps" columns (N, O) pulled in from the
#    weight_df / exercise_df sheets (start of wiring the Postgres-sourced
#    data into the combined calorie dataframe).
# ---------------------------------------------------------------------------
$wsCalorie.Range("N1").Value = "Weight"
$wsCalorie.Range("O1").Value = "Steps"

$wsCalorie.Range("N2").Value = 162.4
$wsCalorie.Range("O2").Value = 12863

$wsCalorie.Range("N3").Value = 161.8
$wsCalorie.Range("O3").Value = 13593

$wsCalorie.Range("N4").Value = 161.4
$wsCalorie.Range("O4").Value = 12128

$wsCalorie.Range("N5").Value = 161
$wsCalorie.Range("O5").Value = 11987

# Resize the columns to fit their (now wider/narrower) contents.
$wsCalorie.Range("A1:O5").EntireColumn.AutoFit()

# ---------------------------------------------------------------------------
# 2. Page setup: force portrait orientation on weight_df / exercise_df.
# ---------------------------------------------------------------------------
$wsWeight.PageSetup.Orientation = 1
$wsExercise.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 3. Restore/update selections on each sheet, finishing on weight_df so it
#    becomes the active tab.
# ---------------------------------------------------------------------------
$wsExercise.Range("A2:B5").Select()
$wsCalorie.Range("E11").Select()
$wsWeight.Range("F24").Select()
